$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that were dropped from the taxonomy table (refactoring for
# reproducibility of the HT analysis paper): the Chiropteran2/DrAAV row, the
# Reptile1/Reptile2 rows, and the two Galliformes1 rows. Delete from the
# bottom up so earlier row numbers stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(3).Delete()
